# Refresh the cryptocurrency price/volume snapshot (and the small block of
# re-ranked coins in rows 16-23) to match the latest scrape.
#
# Numeric-looking strings (prices, percentages) are written with a leading
# apostrophe so Excel keeps them as literal text (matching column D/E's
# existing text formatting) instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''305.18'
$ws.Range("E2").Value = '''0.31%'
# Row 3
$ws.Range("D3").Value = '''37.05'
$ws.Range("E3").Value = '''5.84%'
# Row 4
$ws.Range("D4").Value = '''5.008'
$ws.Range("E4").Value = '''-3.16%'
# Row 5
$ws.Range("D5").Value = '''0.07887'
# Row 6
$ws.Range("D6").Value = '''2.211'
$ws.Range("E6").Value = '''-3.97%'
# Row 7
$ws.Range("D7").Value = '''8.012'
$ws.Range("E7").Value = '''-0.64%'
# Row 8
$ws.Range("D8").Value = '''4.016'
$ws.Range("E8").Value = '''0.71%'
# Row 9
$ws.Range("D9").Value = '''0.9200'
$ws.Range("E9").Value = '''-0.38%'
# Row 10
$ws.Range("D10").Value = '''0.09650'
$ws.Range("E10").Value = '''-4.21%'
# Row 11
$ws.Range("E11").Value = '''3.24%'
# Row 12
$ws.Range("D12").Value = '''0.08595'
$ws.Range("E12").Value = '''0.69%'
# Row 13
$ws.Range("D13").Value = '''0.03681'
$ws.Range("E13").Value = '''8.53%'
# Row 14
$ws.Range("D14").Value = '''0.09986'
$ws.Range("E14").Value = '''0.86%'
# Row 15
$ws.Range("D15").Value = '''0.001479'
$ws.Range("E15").Value = '''-1.28%'
# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005641'
$ws.Range("E16").Value = '''-2.58%'
# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.463'
$ws.Range("E17").Value = '''-0.27%'
# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.249'
$ws.Range("E18").Value = '''7.01%'
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3414'
$ws.Range("E19").Value = '''-0.07%'
# Row 20
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '''0.1316'
$ws.Range("E20").Value = '''-0.77%'
# Row 21
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''4.754'
$ws.Range("E21").Value = '''4.60%'
# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '''0.2198'
$ws.Range("E22").Value = '''-3.21%'
# Row 23
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '''0.04544'
$ws.Range("E23").Value = '''-2.29%'
# Row 24
$ws.Range("D24").Value = '''0.001233'
$ws.Range("E24").Value = '''1.11%'
# Row 25
$ws.Range("D25").Value = '''0.004473'
$ws.Range("E25").Value = '''3.04%'
# Row 26
$ws.Range("E26").Value = '''7.58%'
# Row 27
$ws.Range("E27").Value = '''39.74%'
# Row 39
$ws.Range("D39").Value = '''0.01848'
$ws.Range("E39").Value = '''5.61%'
# Row 40
$ws.Range("D40").Value = '''0.04763'
$ws.Range("E40").Value = '''0.36%'
# Row 41
$ws.Range("D41").Value = '''0.008121'
$ws.Range("E41").Value = '''5.72%'
# Row 42
$ws.Range("D42").Value = '''0.1400'
$ws.Range("E42").Value = '''-0.77%'
# Row 43
$ws.Range("D43").Value = '''0.007555'
$ws.Range("E43").Value = '''-1.24%'
# Row 44
$ws.Range("E44").Value = '''-2.72%'
# Row 45
$ws.Range("D45").Value = '''0.01055'
$ws.Range("E45").Value = '''5.91%'
# Row 46
$ws.Range("D46").Value = '''0.00006291'
$ws.Range("E46").Value = '''3.92%'
# Row 47
$ws.Range("E47").Value = '''-0.13%'
# Row 48
$ws.Range("D48").Value = '''0.0005795'
$ws.Range("E48").Value = '''-0.09%'
# Row 49
$ws.Range("D49").Value = '''29.98'
$ws.Range("E49").Value = '''672.58%'
# Row 50
$ws.Range("D50").Value = '''0.001720'
$ws.Range("E50").Value = '''-36.10%'
# Row 51
$ws.Range("E51").Value = '''-0.13%'
